# Marks sheet correction: fill in missing "Paper 1" / "Paper 2" marks that
# were previously left blank for a number of students on both the
# "Senior Five" and "Senior Six" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "Senior Five" sheet
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Senior Five")

# Column E = "Paper 1" marks that were missing.
$paper1Cells  = @("E2", "E3", "E5", "E6", "E7", "E9", "E11", "E12", "E14", "E15")
$paper1Values = @(42,   56,   62,   30,   58,   42,   52,    56,    38,    62)

for ($i = 0; $i -lt $paper1Cells.Length; $i++) {
    $cell = $ws1.Range($paper1Cells[$i])
    $cell.Value = $paper1Values[$i]
    # Touch the alignment so the cell is normalized onto its own style,
    # distinguishing the now-populated "Paper 1" marks cells.
    $cell.Orientation = 0
}

# Column F = "Paper 2" marks that were missing for a couple of students.
$ws1.Range("F2").Value = 78
$ws1.Range("F7").Value = 80

# ---------------------------------------------------------------------
# "Senior Six" sheet
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Senior Six")

# Column F = "Paper 2" marks that were missing.
$ws2.Range("F3").Value  = 70
$ws2.Range("F5").Value  = 88
$ws2.Range("F6").Value  = 80
$ws2.Range("F8").Value  = 90
$ws2.Range("F9").Value  = 86
$ws2.Range("F12").Value = 70
$ws2.Range("F19").Value = 82
